$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Damian Lillard"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Milwaukee Bucks"

$ws.Range("A3").Value = "Cade Cunningham"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Detroit Pistons"

$ws.Range("A4").Value = "Carlton Carrington"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Washington Wizards"

$ws.Range("A5").Value = "Malik Monk"
$ws.Range("B5").Value = "PG,SG,SF"
$ws.Range("C5").Value = "Sacramento Kings"

$ws.Range("A6").Value = "Devin Vassell"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "San Antonio Spurs"

$ws.Range("A7").Value = "Ausar Thompson"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Detroit Pistons"

$ws.Range("A8").Value = "Santi Aldama"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Memphis Grizzlies"

$ws.Range("A9").Value = "Precious Achiuwa"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "New York Knicks"

$ws.Range("A10").Value = "Isaiah Hartenstein"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Oklahoma City Thunder"

$ws.Range("A11").Value = "Naz Reid"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Minnesota Timberwolves"

$ws.Range("A12").Value = "Andrew Wiggins"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Miami Heat"

$ws.Range("A13").Value = "Derrick White"
$ws.Range("B13").Value = "PG,SG"
$ws.Range("C13").Value = "Boston Celtics"

$ws.Range("A14").Value = "Jusuf Nurkic"
$ws.Range("B14").Value = "C"
$ws.Range("C14").Value = "Charlotte Hornets"

$ws.Range("A15").Value = "Onyeka Okongwu"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Atlanta Hawks"

$ws.Range("A16").Value = "Coby White"
$ws.Range("B16").Value = "PG,SG"
$ws.Range("C16").Value = "Chicago Bulls"

$ws.Range("A17").Value = "Collin Sexton"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Utah Jazz"

$ws.Range("A18").Value = "Anthony Davis"
$ws.Range("B18").Value = "PF,C"
$ws.Range("C18").Value = "Dallas Mavericks"

$ws.Range("A19").Value = "LaMelo Ball"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Charlotte Hornets"
